$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at 13 to shift rows 13-23 down to 14-24 ---
$ws.Rows(13).Insert()

# --- Row 13 (new): only B13/C13 populated, matching style of row 10 (B=s2,C=s3) ---
$ws.Range("A13").Clear()
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update cell text content ---
$txtObjetivos = @"
Desenvolver a capacidade de elaborar, interpretar e executar processos de usinagem. Apresentar e discutir conceitos sobre os processos de usinagem e sobre as interações devido ao contato ferramenta-peça durante as operações de corte. Apresentar os mecanismos de desgaste e avaria além da correlação com a vida da ferramenta. Propor abordagens com enfoque nos aspectos econômicos da usinagem dos materiais. Discutir sobre as novas tendências da tecnologia da usinagem
"@
$ws.Range("B10").Value = $txtObjetivos
$ws.Range("C10").Value = $txtObjetivos

$txtDocente = @"
471420 - Carlos Antonio Reis Pereira Baptista
"@
$ws.Range("B13").Value = $txtDocente
$ws.Range("C13").Value = $txtDocente

$txtResumido = @"
Usinagem: tipos, propriedades, processamento e aplicações. Tendências da tecnologia da usinagem.
"@
$ws.Range("B14").Value = $txtResumido
$ws.Range("C14").Value = $txtResumido

$txtPrograma = @"
Conteúdo teórico: 
1. Processos mecânicos de usinagem. 
2. Mecanismos de formação do cavaco. 
3. Cálculo de potência de corte. 
4. Condições econômicas de corte.
5. Fluido de corte.
6. Usinabilidade dos materiais.
7. Novas tendências na usinagem dos materiais. 
Conteúdo prático: 
1. Trabalhos práticos em máquinas e equipamentos de usinagem.
2. Visita a empresa de usinagem.
"@
$ws.Range("B16").Value = $txtPrograma
$ws.Range("C16").Value = $txtPrograma

$txtMetodo = @"
Aula expositiva com utilização de recursos audiovisuais, aliada a aulas práticas de preparação de máquinas e equipamentos de usinagem e visita a empresas.
"@
$ws.Range("B19").Value = $txtMetodo
$ws.Range("C19").Value = $txtMetodo

$txtCriterio = @"
A nota final será calculada pela média ponderada de duas provas, valendo 60% e da média de exercícios, testes práticos e relatórios de laboratório, valendo 40% da nota final.
A fórmula para o cálculo da média será: NF = (P1+2*P2)/3*0,6 + ME*0,4, na qual P1 e P2 são as notas das provas e ME a média dos exercícios, testes e relatórios.
"@
$ws.Range("B20").Value = $txtCriterio
$ws.Range("C20").Value = $txtCriterio

$txtRecuperacao = @"
A recuperação será uma prova dissertativa. A média final será MF = (NF + RE)/2, na qual NF é a nota final e RE a nota da prova de recuperação.
"@
$ws.Range("B21").Value = $txtRecuperacao
$ws.Range("C21").Value = $txtRecuperacao

$txtBibliografia = @"
1.SCHNEIDER, JR. G., Cutting Tool Applications. Nelson Publishing, Inc. New York, USA, 2001. 
2.DINIZ, A.E.; MARCONDES, F.C.; COPPINI, N.L., Tecnologia da Usinagem dos Materiais. Ed. Artliber, São Paulo, 2ª ed., 2000.
3.FERRAREZI, Dino. Fundamentos da usinagem dos metais, Edgar Blucher, 1995.
4.ABNT - Normas Técnicas de 1995 - edição ABNT.
5.SANDVIK COROMANT, Modern Metal Cutting. AB Sandvik Coromant, Sandviken, Sweden, 1994. 
6.MACHADO, A. Usinagem dos metais. Uberlândia: Universidade Federal de Uberlândia, 1994.
7.DeVRIES, W.R., Analysis of Material Removal Processes. Springer-Verlag, New York, USA, 1991.
"@
$ws.Range("B22").Value = $txtBibliografia
$ws.Range("C22").Value = $txtBibliografia

# --- Fix column definitions: split column A off from the merged A:B width range ---
$ws.Columns("B:B").ColumnWidth = 59.83
